$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.772.42"
$ws.Range("E2").Value = "  -2.30%  "

$ws.Range("D3").Value = "1.750.92"
$ws.Range("E3").Value = "  -4.47%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'236.75"
$ws.Range("E5").Value = "  -6.53%  "

$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").Value = "'0.5075"
$ws.Range("E7").Value = "  -3.92%  "

$ws.Range("D8").Value = "'41.82"
$ws.Range("E8").Value = "  -5.78%  "

$ws.Range("D9").Value = "'0.2642"
$ws.Range("E9").Value = "  -6.78%  "

$ws.Range("D10").Value = "'0.06151"
$ws.Range("E10").Value = "  -10.70%  "

$ws.Range("D11").Value = "1.753.36"
$ws.Range("E11").Value = "  -4.45%  "

$ws.Range("D12").Value = "'15.71"
$ws.Range("E12").Value = "  -5.00%  "

$ws.Range("D13").Value = "'0.06908"
$ws.Range("E13").Value = "  -3.17%  "

$ws.Range("D14").Value = "'0.6052"
$ws.Range("E14").Value = "  -13.14%  "

$ws.Range("D15").Value = "'4.498"
$ws.Range("E15").Value = "  -7.85%  "

$ws.Range("D16").Value = "'77.00"
$ws.Range("E16").Value = "  -11.56%  "

$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  -0.12%  "

$ws.Range("D19").Value = "25.778.75"
$ws.Range("E19").Value = "  -2.35%  "

$ws.Range("D20").Value = "'11.71"
$ws.Range("E20").Value = "  -11.86%  "

$ws.Range("D21").Value = "'0.000006842"
$ws.Range("E21").Value = "  -7.25%  "

$ws.Range("D22").Value = "1.973.34"
$ws.Range("E22").Value = "  -5.39%  "

$ws.Range("D23").Value = "'4.090"
$ws.Range("E23").Value = "  -9.34%  "

$ws.Range("D24").Value = "'8.239"
$ws.Range("E24").Value = "  -8.34%  "

$ws.Range("D25").Value = "'5.210"
$ws.Range("E25").Value = "  -10.63%  "

$ws.Range("D26").Value = "'137.68"
$ws.Range("E26").Value = "  -3.17%  "

$ws.Range("D27").Value = "'1.470"
$ws.Range("E27").Value = "  -12.56%  "

$ws.Range("D28").Value = "'1.821"
$ws.Range("E28").Value = "  -10.84%  "

$ws.Range("D29").Value = "'15.01"
$ws.Range("E29").Value = "  -9.64%  "

$ws.Range("D30").Value = "'102.70"
$ws.Range("E30").Value = "  -5.88%  "

$ws.Range("D31").Value = "'0.08212"
$ws.Range("E31").Value = "  -6.05%  "

$ws.Range("D32").Value = "'3.687"
$ws.Range("E32").Value = "  -9.97%  "

$ws.Range("D33").Value = "'3.462"
$ws.Range("E33").Value = "  -10.57%  "

$ws.Range("E34").Value = "  -4.30%  "

$ws.Range("D35").Value = "'1.000"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").Value = "'2.667"
$ws.Range("E36").Value = "  -7.62%  "

$ws.Range("D37").Value = "'0.9964"
$ws.Range("E37").Value = "  -10.31%  "

$ws.Range("D38").Value = "'0.6071"
$ws.Range("E38").Value = "  -14.18%  "

$ws.Range("D39").Value = "'2.700"
$ws.Range("E39").Value = "  -11.80%  "

$ws.Range("D40").Value = "'0.01554"
$ws.Range("E40").Value = "  -5.73%  "

$ws.Range("D41").Value = "'1.935"
$ws.Range("E41").Value = "  -11.38%  "

$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").Value = "'103.36"
$ws.Range("E43").Value = "  -1.65%  "

$ws.Range("D44").Value = "'0.3821"
$ws.Range("E44").Value = "  -14.77%  "

$ws.Range("D45").Value = "'0.7391"
$ws.Range("E45").Value = "  -14.43%  "

$ws.Range("D46").Value = "'4.926"
$ws.Range("E46").Value = "  -14.46%  "

$ws.Range("D47").Value = "'0.05466"
$ws.Range("E47").Value = "  -2.01%  "

$ws.Range("D48").Value = "'0.1103"
$ws.Range("E48").Value = "  -7.46%  "

$ws.Range("D49").Value = "'5.960"
$ws.Range("E49").Value = "  -15.54%  "

$ws.Range("D50").Value = "'7.675"
$ws.Range("E50").Value = "  -11.36%  "

$ws.Range("D51").Value = "'29.95"
$ws.Range("E51").Value = "  -10.64%  "
